$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.499.76"
$ws.Range("E2").Value = "  +2.11%  "
$ws.Range("D3").Value = "2.324.02"
$ws.Range("E3").Value = "  +1.25%  "
$ws.Range("E4").Value = "  +0.03%  "
$__style = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "515.21"
$ws.Range("D5").Style = $__style
$ws.Range("E5").Value = "  +1.03%  "
$__style = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "134.20"
$ws.Range("D6").Style = $__style
$ws.Range("E6").Value = "  +3.34%  "
$ws.Range("E7").Value = "  +0.33%  "
$ws.Range("E8").Value = "  +1.03%  "
$ws.Range("D9").Value = "2.328.39"
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  +1.85%  "
$__style = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.152"
$ws.Range("D11").Style = $__style
$ws.Range("E11").Value = "  -1.38%  "
$__style = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.35"
$ws.Range("D12").Style = $__style
$ws.Range("E12").Value = "  +4.80%  "
$ws.Range("E13").Value = "  -0.16%  "
$__style = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.86"
$ws.Range("D14").Style = $__style
$ws.Range("E14").Value = "  +0.33%  "
$ws.Range("D15").Value = "2.731.82"
$ws.Range("E15").Value = "  +1.55%  "
$ws.Range("D16").Value = "56.553.95"
$ws.Range("E16").Value = "  +2.22%  "
$ws.Range("E17").Value = "  +1.33%  "
$ws.Range("D18").Value = "2.330.32"
$ws.Range("E18").Value = "  +2.16%  "
$ws.Range("E19").Value = "  +0.26%  "
$__style = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "325.66"
$ws.Range("D20").Style = $__style
$ws.Range("E20").Value = "  +2.73%  "
$ws.Range("E21").Value = "  +0.04%  "
$__style = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.59"
$ws.Range("D22").Style = $__style
$ws.Range("E22").Value = "  -0.02%  "
$__style = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.999"
$ws.Range("D23").Style = $__style
$ws.Range("E23").Value = "  +0.07%  "
$__style = $ws.Range("D24").Style
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.66"
$ws.Range("D24").Style = $__style
$ws.Range("E24").Value = "  +1.06%  "
$__style = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.164"
$ws.Range("D25").Style = $__style
$ws.Range("E25").Value = "  +5.13%  "
$ws.Range("E26").Value = "  +0.20%  "
$__style = $ws.Range("D27").Style
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "8.06"
$ws.Range("D27").Style = $__style
$ws.Range("E27").Value = "  +6.10%  "
$__style = $ws.Range("D28").Style
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.29"
$ws.Range("D28").Style = $__style
$ws.Range("E28").Value = "  +11.22%  "
$__style = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "169.82"
$ws.Range("D29").Style = $__style
$ws.Range("E29").Value = "  -0.76%  "
$ws.Range("D30").Value = "0.0₃0737"
$ws.Range("E30").Value = "  +3.27%  "
$ws.Range("E31").Value = "  +2.11%  "
$__style = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.21"
$ws.Range("D32").Style = $__style
$__style = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "18.48"
$ws.Range("D33").Style = $__style
$ws.Range("E33").Value = "  +2.13%  "
$ws.Range("E34").Value = "  +0.02%  "
$__style = $ws.Range("D35").Style
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.997"
$ws.Range("D35").Style = $__style
$ws.Range("E35").Value = "  +0.84%  "
$__style = $ws.Range("D36").Style
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.26"
$ws.Range("D36").Style = $__style
$ws.Range("E36").Value = "  +0.53%  "
$ws.Range("E37").Value = "  -0.07%  "
$ws.Range("E38").Value = "  +1.73%  "
$ws.Range("E39").Value = "  +5.25%  "
$__style = $ws.Range("D40").Style
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "38.12"
$ws.Range("D40").Style = $__style
$ws.Range("E40").Value = "  +3.21%  "
$ws.Range("E41").Value = "  +0.97%  "
$__style = $ws.Range("D42").Style
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "140.51"
$ws.Range("D42").Style = $__style
$ws.Range("E42").Value = "  +2.04%  "
$__style = $ws.Range("D43").Style
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.60"
$ws.Range("D43").Style = $__style
$ws.Range("E43").Value = "  +2.53%  "
$ws.Range("B44").Value = "Bittensor"
$ws.Range("C44").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$__style = $ws.Range("D44").Style
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "277.84"
$ws.Range("D44").Style = $__style
$ws.Range("E44").Value = "  +5.82%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$__style = $ws.Range("D45").Style
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "5.13"
$ws.Range("D45").Style = $__style
$ws.Range("E45").Value = "  +1.22%  "
$__style = $ws.Range("D46").Style
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0935"
$ws.Range("D46").Style = $__style
$ws.Range("E47").Value = "  -0.72%  "
$ws.Range("E48").Value = "  +1.43%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$__style = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "17.78"
$ws.Range("D49").Style = $__style
$ws.Range("E49").Value = "  +6.61%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$__style = $ws.Range("D50").Style
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0218"
$ws.Range("D50").Style = $__style
$ws.Range("E50").Value = "  +1.63%  "
$__style = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "17.83"
$ws.Range("D51").Style = $__style
$ws.Range("E51").Value = "  +7.10%  "
